$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Era C)
$ws.Range("B2").Value = 346021.01
$ws.Range("C2").Value = 340067.57
$ws.Range("G2").Value = 686088.58

# Row 3 (Era B)
$ws.Range("D3").Value = 406645.98
$ws.Range("G3").Value = 406645.98

# Row 4 (Era A)
$ws.Range("E4").Value = 271806.75
$ws.Range("F4").Value = 426460.18
$ws.Range("G4").Value = 698266.9300000001

# Row 5 (Total)
$ws.Range("B5").Value = 346021.01
$ws.Range("C5").Value = 340067.57
$ws.Range("D5").Value = 406645.98
$ws.Range("E5").Value = 271806.75
$ws.Range("F5").Value = 426460.18
$ws.Range("G5").Value = 1791001.49
